$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort/update the weekly price data rows (2-23) to match the new snapshot.
# Each row is fully rewritten (date, variety, quality, volume, prices, unit,
# origin, $/Kg, Kg/unit) since the underlying weekly extract re-ordered rows.

# Row 2
$ws.Range("D2").Value = 44497
$ws.Range("K2").Value = 'Packham''s Triumph'
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 17500
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 972
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44497
$ws.Range("K3").Value = 'Winter Nelis'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17500
$ws.Range("Q3").Value = '$/bandeja 18 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 972
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44355
$ws.Range("K4").Value = 'Packham''s Triumph'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 972
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44355
$ws.Range("K5").Value = 'Winter Nelis'
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17500
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 972
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44313
$ws.Range("K6").Value = 'Winter Nelis'
$ws.Range("L6").Value = 'Tercera'
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44336
$ws.Range("K7").Value = 'Winter Nelis'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1194
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44314
$ws.Range("K8").Value = 'Packham''s Triumph'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 972
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44474
$ws.Range("K9").Value = 'Packham''s Triumph'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 270
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("Q9").Value = '$/caja 18 kilos empedrada'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1028
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44474
$ws.Range("K10").Value = 'Winter Nelis'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 17500
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 972
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44341
$ws.Range("K11").Value = 'Packham''s Triumph'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 17000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 17500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 972
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44398
$ws.Range("K12").Value = 'Packham''s Triumph'
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("Q12").Value = '$/caja 20 kilos empedrada'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 1025
$ws.Range("T12").Value = 20

# Row 13
$ws.Range("D13").Value = 44398
$ws.Range("K13").Value = 'Winter Nelis'
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("Q13").Value = '$/caja 20 kilos empedrada'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1025
$ws.Range("T13").Value = 20

# Row 14
$ws.Range("D14").Value = 44371
$ws.Range("K14").Value = 'Packham''s Triumph'
$ws.Range("L14").Value = 'Calibre 90'
$ws.Range("M14").Value = 140
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17429
$ws.Range("Q14").Value = '$/caja 18 kilos embalada'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 968
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("D15").Value = 44371
$ws.Range("K15").Value = 'Winter Nelis'
$ws.Range("L15").Value = 'Calibre 80'
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 17500
$ws.Range("Q15").Value = '$/caja 18 kilos embalada'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 972
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("D16").Value = 44329
$ws.Range("K16").Value = 'Packham''s Triumph'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 340
$ws.Range("N16").Value = 21000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 21500
$ws.Range("Q16").Value = '$/bandeja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 1194
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44280
$ws.Range("K17").Value = 'Packham''s Triumph'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 350
$ws.Range("N17").Value = 24000
$ws.Range("O17").Value = 25000
$ws.Range("P17").Value = 24500
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1361
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = 44280
$ws.Range("K18").Value = 'Winter Nelis'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 25000
$ws.Range("P18").Value = 24500
$ws.Range("Q18").Value = '$/caja 18 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1361
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44292
$ws.Range("K19").Value = 'Packham''s Triumph'
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 22000
$ws.Range("O19").Value = 23000
$ws.Range("P19").Value = 22500
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 1250
$ws.Range("T19").Value = 18

# Row 20
$ws.Range("D20").Value = 44292
$ws.Range("K20").Value = 'Winter Nelis'
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 250
$ws.Range("N20").Value = 22000
$ws.Range("O20").Value = 23000
$ws.Range("P20").Value = 22500
$ws.Range("Q20").Value = '$/caja 18 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 1250
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("D21").Value = 44323
$ws.Range("K21").Value = 'Packham''s Triumph'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 250
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 16000
$ws.Range("P21").Value = 15500
$ws.Range("Q21").Value = '$/bandeja 18 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 861
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("D22").Value = 44421
$ws.Range("K22").Value = 'Packham''s Triumph'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 270
$ws.Range("N22").Value = 16000
$ws.Range("O22").Value = 17000
$ws.Range("P22").Value = 16500
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 917
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("D23").Value = 44421
$ws.Range("K23").Value = 'Winter Nelis'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 16000
$ws.Range("O23").Value = 17000
$ws.Range("P23").Value = 16500
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 917
$ws.Range("T23").Value = 18
